$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.146576046943665
$ws.Range("B1").Value = 3.463815450668335
$ws.Range("C1").Value = 3.452594518661499
$ws.Range("D1").Value = 3.860114812850952
$ws.Range("E1").Value = 1.130586504936218
